# Revert config file handling
# Appends one new trailing row (row 45) to each of the four worksheets,
# duplicating the prior last row (row 44) but with an updated timestamp
# in column A (one day's worth of new log data appended).

$wb = $excel.ActiveWorkbook

# Scientific-notation literals aren't accepted directly by the parser,
# so build the doubles from strings once up front.
$gLft1 = [double]"7.598631275147109e+23"
$gLft2 = [double]"5.68432987514711e+23"
$gPlt1 = [double]"5.68631262647114e+23"
$gPlt2 = [double]"9.85046333984776e+23"

# --- Sheet 1: DE_LFT_#1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A45").NumberFormat = $ws1.Range("A44").NumberFormat
$ws1.Range("A45").Value2 = 45831.43658564815
$ws1.Range("B45").Value = "0x01,0x7c"
$ws1.Range("C45").Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws1.Range("D45").Value = "0x01,0x68"
$ws1.Range("E45").Value = "0x14"
$ws1.Range("F45").Value = 380
$ws1.Range("G45").Value = $gLft1
$ws1.Range("H45").Value = 360
$ws1.Range("I45").Value = 14

# --- Sheet 2: DE_LFT_#2 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A45").NumberFormat = $ws2.Range("A44").NumberFormat
$ws2.Range("A45").Value2 = 45831.43658564815
$ws2.Range("B45").Value = "0x01,0x7c"
$ws2.Range("C45").Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws2.Range("D45").Value = "0x01,0x68"
$ws2.Range("E45").Value = "0xe"
$ws2.Range("F45").Value = 380
$ws2.Range("G45").Value = $gLft2
$ws2.Range("H45").Value = 360
$ws2.Range("I45").Value = 14

# --- Sheet 3: DE_PLT_#1 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A45").NumberFormat = $ws3.Range("A44").NumberFormat
$ws3.Range("A45").Value2 = 45831.43658564815
$ws3.Range("B45").Value = "0x00,0x82"
$ws3.Range("C45").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Range("D45").Value = "0x00,0x7F"
$ws3.Range("E45").Value = "0x7"
$ws3.Range("F45").Value = 130
$ws3.Range("G45").Value = $gPlt1
$ws3.Range("H45").Value = 127
$ws3.Range("I45").Value = 7

# --- Sheet 4: DE_PLT_#2 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A45").NumberFormat = $ws4.Range("A44").NumberFormat
$ws4.Range("A45").Value2 = 45831.43658564815
$ws4.Range("B45").Value = "0x00,0x82"
$ws4.Range("C45").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Range("D45").Value = "0x00,0x7F"
$ws4.Range("E45").Value = "0x3"
$ws4.Range("F45").Value = 130
$ws4.Range("G45").Value = $gPlt2
$ws4.Range("H45").Value = 127
$ws4.Range("I45").Value = 3
